# Update the "F" (sales count) column values on the "展览" and "全部类型"
# sheets to reflect refreshed sales figures.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 593
    $ws.Range("F3").Value = 130

    if ($name -eq "展览") {
        $ws.Range("F7").Value = 1727
        $ws.Range("F8").Value = 97
    }
    else {
        $ws.Range("F11").Value = 1727
        $ws.Range("F12").Value = 97
    }
}
